$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (task descriptions) -----------------------------------------
# Rows 3, 7, 8, 9, 10, 11 keep the exact same text as before, so they are
# left untouched (this also keeps those shared-string entries alive/reused).
#
# The order below matters: new/unique strings are appended to the shared
# string table in first-write order, so we walk the sheet top-to-bottom.

$ws.Cells.Item(4, 2).Value  = "Start to work on Alice"
$ws.Cells.Item(5, 2).Value  = "Have a break, Go for lunch"
$ws.Cells.Item(6, 2).Value  = "Contribution in Python-Discord"

$ws.Cells.Item(12, 2).Value = "Go for the dinner by shutting down pc"
$ws.Cells.Item(13, 2).Value = "Time pass on python discord by helping others"
$ws.Cells.Item(14, 2).Value = "ATL Session"
$ws.Cells.Item(15, 2).Value = "Playing Krunker online Game"
$ws.Cells.Item(16, 2).Value = "Start Solving the Practise sets of python"
$ws.Cells.Item(17, 2).Value = "Read Books"

# --- Column A (time-of-day values) -----------------------------------------
$ws.Cells.Item(12, 1).Value = 0.8125
$ws.Cells.Item(13, 1).Value = 0.71875
$ws.Cells.Item(13, 1).NumberFormat = "h:mm:ss"
$ws.Cells.Item(14, 1).Value = 0.79166666666666663
$ws.Cells.Item(14, 1).NumberFormat = "h:mm:ss"
$ws.Cells.Item(15, 1).Value = 0.67708333333333337
$ws.Cells.Item(16, 1).Value = 0.90625
$ws.Cells.Item(16, 1).NumberFormat = "h:mm AM/PM"
$ws.Cells.Item(17, 1).Value = 0.95833333333333337
$ws.Cells.Item(17, 1).NumberFormat = "h:mm AM/PM"

# --- New rows 18 & 19 --------------------------------------------------------
$ws.Cells.Item(18, 1).Value = 0.35416666666666669
$ws.Cells.Item(18, 1).NumberFormat = "h:mm AM/PM"
$ws.Cells.Item(18, 2).Value = "Read core python book"

$ws.Cells.Item(19, 1).Value = 0.41666666666666669
$ws.Cells.Item(19, 1).NumberFormat = "h:mm"
$ws.Cells.Item(19, 2).Value = "Scroll Instagram and post 1 python trick "

# --- Selection mirrors the post-edit cursor position ------------------------
$ws.Range("B32").Select()
